$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new project row at row 20 by shifting the existing rows 20-23
# down to 21-24 (processed bottom-up so later moves don't clobber data
# that hasn't been relocated yet). Using Cut/Paste (instead of
# Rows.Insert) keeps the original cell styles (s="4"/"8"/"12") intact
# without Excel minting brand-new style records.
$ws.Range("A23:J23").Cut($ws.Range("A24:J24"))
$ws.Range("A22:J22").Cut($ws.Range("A23:J23"))
$ws.Range("A21:J21").Cut($ws.Range("A22:J22"))
$ws.Range("A20:J20").Cut($ws.Range("A21:J21"))

# Populate the freed-up row 20 with the new PROCORTE project entry.
$ws.Range("A20").Value = "procorte"
$ws.Range("B20").Value = "Sistema com algoritmos avançados para otimização de problemas de corte"
$ws.Range("C20").Value = (Get-Date -Year 2017 -Month 7 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D20").Value = "PROT"
$ws.Range("E20").Value = "X"
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = ""
$ws.Range("I20").Value = "MPS et al."
$ws.Range("J20").Value = ""

$ws.Range("J20").Select() | Out-Null
